$wb = $excel.ActiveWorkbook

# Excel's COM ColumnWidth setter/getter has a fixed +5/6 (~0.8333) offset that shows
# up when the width is serialized back to the OOXML "width" attribute. Subtracting
# that offset before assigning lets us land on an exact integer width in the saved file.
$widthOffset = 5/6

function Set-ColWidth($ws, $col, $target) {
    $ws.Columns.Item($col).ColumnWidth = $target - $widthOffset
}

function Add-PriceRow($ws, $row, $dateStr, $price) {
    $ws.Cells.Item($row, 1).Value = $dateStr
    $ws.Cells.Item($row, 2).Value = $price
}

function Add-DateOnlyRow($ws, $row, $dateStr) {
    $ws.Cells.Item($row, 1).Value = $dateStr
}

function Set-ChartRange($ws, $sheetName, $lastRow) {
    $co = $ws.ChartObjects().Item(1)
    $chart = $co.Chart
    $ser = $chart.SeriesCollection().Item(1)
    $ser.Formula = "=SERIES(,'" + $sheetName + "'!`$A`$4:`$A`$" + $lastRow + ",'" + $sheetName + "'!`$B`$4:`$B`$" + $lastRow + ",1)"
}

# ---------------------------------------------------------------------------
# Sheet 1: "18 - inSPORTline Odino"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("18 - inSPORTline Odino")

Set-ColWidth $ws1 5 9

$ws1.Cells.Item(4, 5).Value = 1808.75

Add-PriceRow $ws1 45 "2026-01-18T16:43:18" 1849
Add-PriceRow $ws1 46 "2026-01-18T16:48:43" 1849
Add-PriceRow $ws1 47 "2026-01-18T16:54:10" 1849
Add-PriceRow $ws1 48 "2026-01-18T16:55:50" 1849
Add-PriceRow $ws1 49 "2026-01-18T16:58:52" 1849

Set-ChartRange $ws1 "18 - inSPORTline Odino" 149

# ---------------------------------------------------------------------------
# Sheet 2: "19 - Produkt"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("19 - Produkt")

Add-DateOnlyRow $ws2 40 "2026-01-18T16:43:42"
Add-DateOnlyRow $ws2 41 "2026-01-18T16:49:07"
Add-DateOnlyRow $ws2 42 "2026-01-18T16:54:33"
Add-DateOnlyRow $ws2 43 "2026-01-18T16:56:13"
Add-DateOnlyRow $ws2 44 "2026-01-18T16:59:16"

Set-ChartRange $ws2 "19 - Produkt" 149

# ---------------------------------------------------------------------------
# Sheet 3: "20 - Marshall Major IV Bluetoo"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("20 - Marshall Major IV Bluetoo")

Set-ColWidth $ws3 5 8

$ws3.Cells.Item(4, 5).Value = 1393.6
$ws3.Cells.Item(5, 5).Value = 1399

Add-PriceRow $ws3 35 "2026-01-18T16:44:07" 1399
Add-PriceRow $ws3 36 "2026-01-18T16:49:31" 1399
Add-PriceRow $ws3 37 "2026-01-18T16:54:57" 1399
Add-PriceRow $ws3 38 "2026-01-18T16:56:37" 1399
Add-PriceRow $ws3 39 "2026-01-18T16:59:43" 1399

Set-ChartRange $ws3 "20 - Marshall Major IV Bluetoo" 149

# ---------------------------------------------------------------------------
# Sheet 4: "21 - PlayStation 5 Slim"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("21 - PlayStation 5 Slim")

Set-ColWidth $ws4 5 20

$ws4.Cells.Item(4, 5).Value = 12288.47368421053

Add-PriceRow $ws4 24 "2026-01-18T16:44:16" 12423
Add-PriceRow $ws4 25 "2026-01-18T16:49:49" 12423
Add-PriceRow $ws4 26 "2026-01-18T16:55:10" 12423
Add-PriceRow $ws4 27 "2026-01-18T16:56:46" 12423
Add-PriceRow $ws4 28 "2026-01-18T16:59:52" 12423

Set-ChartRange $ws4 "21 - PlayStation 5 Slim" 149

# ---------------------------------------------------------------------------
# Sheet 5: "22 - Marshall Major IV Bluetoo"
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("22 - Marshall Major IV Bluetoo")

Set-ColWidth $ws5 5 11

$ws5.Cells.Item(4, 5).Value = 1438.4375
$ws5.Cells.Item(5, 5).Value = 1399

Add-PriceRow $ws5 21 "2026-01-18T16:44:24" 1399
Add-PriceRow $ws5 22 "2026-01-18T16:49:58" 1399
Add-PriceRow $ws5 23 "2026-01-18T16:55:19" 1399
Add-PriceRow $ws5 24 "2026-01-18T16:56:54" 1399
Add-PriceRow $ws5 25 "2026-01-18T17:00:01" 1399

Set-ChartRange $ws5 "22 - Marshall Major IV Bluetoo" 149
